$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Result column's type string to include a unit
$ws.Range("I2").Value = "#float,  unit:ng/µl"

# Add a new third row containing enum/description metadata for each column
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"

# H3/I3 hold an empty string - use a quote-prefixed empty value so the
# cell is a real (shared-string) text cell rather than a cleared/blank
# cell, then drop the resulting quote-prefix formatting.
$ws.Cells.Item(3, 8).Value = "'"
$ws.Cells.Item(3, 9).Value = "'"
$ws.Cells.Item(3, 8).ClearFormats()
$ws.Cells.Item(3, 9).ClearFormats()
